$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original Text formatting so
# numeric-looking values (e.g. "19.32") are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.220.29'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '1.588.22'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '212.10'
$ws.Range('E5').Value = '  +1.61%  '
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +0.67%  '
$ws.Range('D9').Value = '0.0607'
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('D10').Value = '19.32'
$ws.Range('E10').Value = '  -1.27%  '
$ws.Range('D11').Value = '0.0848'
$ws.Range('E11').Value = '  +0.48%  '
$ws.Range('D12').Value = '1.810.60'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').Value = '1.593.32'
$ws.Range('E13').Value = '  +1.58%  '
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('D15').Value = '0.519'
$ws.Range('E15').Value = '  +1.02%  '
$ws.Range('D16').Value = '64.23'
$ws.Range('E16').Value = '  -0.16%  '
$ws.Range('D17').Value = '26.232.63'
$ws.Range('E17').Value = '  +0.36%  '
$ws.Range('D18').Value = '0.0₃0726'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('E19').Value = '  +1.45%  '
$ws.Range('D20').Value = '213.07'
$ws.Range('E20').Value = '  +2.61%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = '4.26'
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  +1.86%  '
$ws.Range('D25').Value = '143.51'
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').Value = '7.02'
$ws.Range('E27').Value = '  +0.59%  '
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('D29').Value = '15.17'
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('E30').Value = '  -1.81%  '
$ws.Range('E31').Value = '  +1.31%  '
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('D33').Value = '1.338.76'
$ws.Range('E33').Value = '  +4.85%  '
$ws.Range('E34').Value = '  -1.41%  '
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('E36').Value = '  -0.79%  '
$ws.Range('E37').Value = '  -6.30%  '
$ws.Range('E38').Value = '  +0.59%  '
$ws.Range('D39').Value = '0.824'
$ws.Range('E39').Value = '  +1.87%  '
$ws.Range('D40').Value = '5.75'
$ws.Range('E40').Value = '  +3.29%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').Value = '0.971'
$ws.Range('E42').Value = '  -12.60%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '2.14'
$ws.Range('E43').Value = '  +0.43%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = '0.767'
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('D45').Value = '1.722.80'
$ws.Range('E45').Value = '  +0.80%  '
$ws.Range('D46').Value = '61.10'
$ws.Range('E46').Value = '  -1.88%  '
$ws.Range('D47').Value = '85.70'
$ws.Range('E47').Value = '  -3.32%  '
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('E49').Value = '  -2.60%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.0977'
$ws.Range('E50').Value = '  -2.66%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.0501'
$ws.Range('E51').Value = '  -0.73%  '
